# Update values in column A for specific rows as described in the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -21.12900000000001
$ws.Range("A6").Value = -20.04169999999999
$ws.Range("A7").Value = -21.28230000000001
$ws.Range("A8").Value = -20.34699999999999
$ws.Range("A16").Value = -20.12419999999999
$ws.Range("A20").Value = -22.29730000000003
$ws.Range("A21").Value = -20.39119999999999
